# Update the "Förändrad" (Changed) date column (C) for rows 2-17
# from 2023-09-06 (serial 45175) to 2023-09-08 (serial 45177).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45177
    }
}
